$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header title
$ws.Range("A1").Value = "Commodities Down in Price"

# Update row 2 data and convert B2 to a real number
$ws.Range("A2").Value = "Steel — Carbon*."
$ws.Range("B2").Value = 1

# Delete rows 3 through 13 (old data rows no longer needed)
$ws.Range("A3:B13").EntireRow.Delete()
